$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d10Text = @'
Για τον υπολογιστή **172.16.150.10/20** να υπολογίσετε:  
**Δ1.** Την μάσκα δικτύου(δυαδική-δεκαδική)  
**Δ2.** Τη διεύθυνση δικτύου (network address)   
**Δ3.** Τη διεύθυνση εκπομπής (broadcast address)  
**Δ4.** Τον συνολικό αριθμό υπολογιστών του συγκεκριμένου δικτύου   
**Δ5.** Την περιοχή διευθύνσεων για υπολογιστές (από-έως) οι οποίες ανήκουν στο ίδιο δίκτυο με τον συγκεκριμένο υπολογιστή  
'@

$e10Text = @'
Για τον υπολογιστή **172.16.150.10/20** να υπολογίσετε:  
**Δ1**. Την μάσκα δικτύου(δυαδική-δεκαδική)  
255.255.240.0 ή 11111111.11111111.1111 **0000.00000000**  
**Δ2.** Τη διεύθυνση δικτύου (network address)  
Διεύθυνση Δικτύου   
172.16.150.10(1010 **0110.00001010**) AND 255.255.240.0(1111 **00000.00000000**)=172.16.144.0(1001 **0000.00000000**)/20  
**Δ3.** Τη διεύθυνση εκπομπής (broadcast address)  
Διεύθυνση Εκπομπής ->  172.16.159.255(1001 **1111.11111111**)  
**Δ4.** Τον συνολικό αριθμό υπολογιστών του συγκεκριμένου δικτύου  
2^12 -2 = 4094  
**Δ5.** Την περιοχή διευθύνσεων για υπολογιστές (από-έως) οι οποίες ανήκουν στο ίδιο δίκτυο με τον συγκεκριμένο υπολογιστή  
Από 172.16.144.1 έως 172.16.159.254  
'@

$originalRowHeight = $ws.Rows.Item(10).RowHeight

$ws.Range("D10").Value = $d10Text
$ws.Range("E10").Value = $e10Text

$ws.Rows.Item(10).RowHeight = $originalRowHeight
